$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48 (pushes old rows 48..157 down to 49..158)
$ws.Rows(48).Insert()

# Insert a new row at position 118 (pushes old row 117 -- now sitting at 118 -- down to 119,
# and everything below it down as well)
$ws.Rows(118).Insert()

# Fill in the two brand-new data rows with the same constant columns as every
# other row in this sheet, plus their own D/J/K/L/M/P values.

# New row 48
$ws.Range("A48").Value = 5
$ws.Range("B48").Value = 'Macroferia Regional de Talca'
$ws.Range("C48").Value = 'Maule'
$ws.Range("D48").Value = 44925
$ws.Range("E48").Value = 7
$ws.Range("F48").Value = 100112030
$ws.Range("G48").Value = 'Poroto granado'
$ws.Range("H48").Value = 'Sin especificar'
$ws.Range("I48").Value = 'Primera'
$ws.Range("J48").Value = 500
$ws.Range("K48").Value = 30000
$ws.Range("L48").Value = 30000
$ws.Range("M48").Value = 30000
$ws.Range("N48").Value = '$/saco 25 kilos'
$ws.Range("O48").Value = 'Región del Maule'
$ws.Range("P48").Value = 1200
$ws.Range("Q48").Value = 25
$ws.Range("R48").Value = 'Hortaliza'

# New row 118
$ws.Range("A118").Value = 5
$ws.Range("B118").Value = 'Macroferia Regional de Talca'
$ws.Range("C118").Value = 'Maule'
$ws.Range("D118").Value = 44924
$ws.Range("E118").Value = 7
$ws.Range("F118").Value = 100112030
$ws.Range("G118").Value = 'Poroto granado'
$ws.Range("H118").Value = 'Sin especificar'
$ws.Range("I118").Value = 'Primera'
$ws.Range("J118").Value = 300
$ws.Range("K118").Value = 30000
$ws.Range("L118").Value = 30000
$ws.Range("M118").Value = 30000
$ws.Range("N118").Value = '$/saco 25 kilos'
$ws.Range("O118").Value = 'Región del Maule'
$ws.Range("P118").Value = 1200
$ws.Range("Q118").Value = 25
$ws.Range("R118").Value = 'Hortaliza'
